$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet by cloning "2021-Q4" (so it starts
#    out with identical layout/number-formats/column structure) and
#    dropping it in right after "2021-Q4" (i.e. before "总计").
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# The clone has 4 data rows (rows 2-5); 2022-Q1 only needs 3 (rows 2-4)
$new.Rows.Item(5).Delete()

# ------------------------------------------------------------------
# 2. Fill in the 2022-Q1 fund holdings.
#    Columns B-G hold text (fund code/name/figures as strings), so we
#    temporarily mark the range as Text before writing, then clear the
#    format again afterwards so the cells end up unstyled, exactly
#    like the source sheet's body cells.
# ------------------------------------------------------------------
$new.Range("B2:G4").NumberFormat = "@"

$new.Range("B2").Value = "013067"
$new.Range("C2").Value = "富安达中小盘六个月持有期混合"
$new.Range("D2").Value = "2.45"
$new.Range("E2").Value = "74.39"
$new.Range("F2").Value = "2.75"
$new.Range("G2").Value = "0.0674"

$new.Range("B3").Value = "008422"
$new.Range("C3").Value = "中融研发创新混合A"
$new.Range("D3").Value = "2.32"
$new.Range("E3").Value = "32.66"
$new.Range("F3").Value = "1.98"
$new.Range("G3").Value = "0.0459"

$new.Range("B4").Value = "008423"
$new.Range("C4").Value = "中融研发创新混合C"
$new.Range("D4").Value = "0.62"
$new.Range("E4").Value = "32.66"
$new.Range("F4").Value = "1.98"
$new.Range("G4").Value = "0.0123"

$new.Range("B2:G4").ClearFormats()

# Column A (row index) and H (position rank) are plain numbers.
$new.Range("A2").Value = 0
$new.Range("A3").Value = 1
$new.Range("A4").Value = 2

$new.Range("H2").Value = 8
$new.Range("H3").Value = 8
$new.Range("H4").Value = 8

# ------------------------------------------------------------------
# 3. Update the "总计" (totals) sheet: insert a new top data row for
#    2022-Q1, pushing the existing 2021-Q4 / 2021-Q3 rows down.
# ------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()
$tot.Range("A2:D2").ClearFormats()

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 3
$tot.Range("D2").Value = 0.13

$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2

# Re-apply the row-label style (s="2") to A2 to match A3/A4, by
# copying the format from the cell below.
$tot.Range("A4").Copy()
$tot.Range("A2").PasteSpecial(-4122)
$tot.Range("A2").Value = 0
